# -----------------------------------------------------------------------
# Update the beta-feedback paragraph:
#   "The climate clause picker has been launched in beta. To help us
#    improve this tool and develop it further, we would be grateful for
#    your feedback. Please contact "
# becomes
#   "The climate clause picker has been launched in beta." + [_GoBack] +
#   " It collates into one downloadable document, TCLP clauses published
#    on our website as at 7 April 2022. Clauses published and updated
#    since 7 April 2022 are not yet included in the Climate Clause
#    Selector. We are working on adding them to the Climate Clause
#    Selector. In the meantime, to help us improve this tool and develop
#    it further, we would be grateful for your feedback. Please contact "
# and Word's "last edit" (_GoBack) bookmark is relocated here. Moving
# the bookmark also renumbers the remaining bookmark ids up by one
# (matching the renumbered _Toc90128866 / _Toc90128590 bookmarks later
# in the document) and leaves the paragraph that used to host _GoBack
# empty.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Step 1 - insert the new explanatory sentence plus the comma that turns
# the old "To" into ", to", right before the existing "To help us
# improve..." text (which is left otherwise untouched).
$insRng = $d.Content
$insRng.Find.Execute("To help us improve")
$insPoint = $d.Range($insRng.Start, $insRng.Start)
$insPoint.InsertBefore(
    "It collates into one downloadable document, TCLP clauses published " + `
    "on our website as at 7 April 2022. Clauses published and updated " + `
    "since 7 April 2022 are not yet included in the Climate Clause " + `
    "Selector. We are working on adding them to the Climate Clause " + `
    "Selector. In the meantime, t")

# Step 2 - the original text started with a capital "To"; drop the old
# capital "T" now that ", t" (lower-case) precedes "o help us improve...".
$tRng = $d.Content
$tRng.Find.Execute("meantime, t")
$oldT = $d.Range($tRng.End, $tRng.End + 1)
$oldT.Delete()

# Step 3 - relocate the "_GoBack" bookmark to sit right after "beta."
# (immediately before the space that starts the new sentence). Adding a
# bookmark whose name already exists elsewhere moves it here and
# renumbers bookmark ids accordingly.
$bmRng = $d.Content
$bmRng.Find.Execute("beta.")
$goBackPoint = $d.Range($bmRng.End, $bmRng.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# Step 4 - splitting the run at the bookmark leaves a stray
# xml:space="preserve" on "The climate clause picker has been launched
# in beta." even though it has no leading/trailing whitespace; do a
# genuine round-trip replace (through a scratch value) so the run's
# text gets freshly re-serialized without the redundant attribute.
$d.Content.Find.Execute(
    "The climate clause picker has been launched in beta.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "TCLP_TMP_MARKER", 2)
$d.Content.Find.Execute(
    "TCLP_TMP_MARKER",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The climate clause picker has been launched in beta.", 2)
